# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.579.81"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "3.367.94"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.07"
$ws.Range("E5").Value = "  -2.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.12"
$ws.Range("E6").Value = "  +0.82%  "

$ws.Range("E7").Value = "  -0.59%  "

$ws.Range("D8").Value = "3.358.52"
$ws.Range("E8").Value = "  -2.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("E11").Value = "  +3.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.70"
$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.09"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").Value = "3.897.83"
$ws.Range("E15").Value = "  -2.62%  "

$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("E17").Value = "  -1.87%  "

$ws.Range("D18").Value = "3.353.33"
$ws.Range("E18").Value = "  -3.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.90"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").Value = "64.513.86"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.986"
$ws.Range("E21").Value = "  -0.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "462.88"
$ws.Range("E22").Value = "  +13.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.87"
$ws.Range("E23").Value = "  +9.68%  "

$ws.Range("E24").Value = "  -2.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.56"
$ws.Range("E25").Value = "  +2.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.45"
$ws.Range("E26").Value = "  +0.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.86"
$ws.Range("E28").Value = "  +1.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.83"
$ws.Range("E29").Value = "  -2.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.12"
$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.53"
$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "583.10"
$ws.Range("E33").Value = "  -0.56%  "

$ws.Range("E34").Value = "  -0.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.69"
$ws.Range("E35").Value = "  -1.70%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("E37").Value = "  -8.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.48"
$ws.Range("E38").Value = "  -1.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.85"
$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("D40").Value = "0.0₃0758"
$ws.Range("E40").Value = "  -2.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.371"
$ws.Range("E41").Value = "  -1.44%  "

$ws.Range("D42").Value = "3.109.06"
$ws.Range("E42").Value = "  -2.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -3.89%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0411"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.22"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("E49").Value = "  -1.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.33"
$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.32"
$ws.Range("E51").Value = "  -1.04%  "
